$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Single quote character used as Excel text-qualifier prefix so that
# numeric-looking "Price" strings (e.g. "1.102", "0.9993") are stored
# literally as text, matching the source data, instead of being
# auto-converted to numbers by Excel.
$q = "'"

$ws.Range('D2').Value = $q + '28.104.87'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = $q + '1.813.29'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').Value = $q + '0.9993'
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').Value = $q + '310.09'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').Value = $q + '1.000'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').Value = $q + '0.5103'
$ws.Range('E7').Value = '  -3.15%  '
$ws.Range('D8').Value = $q + '0.3907'
$ws.Range('E8').Value = '  +2.15%  '
$ws.Range('D9').Value = $q + '0.09778'
$ws.Range('E9').Value = '  +21.94%  '
$ws.Range('D10').Value = $q + '1.102'
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('D11').Value = $q + '40.96'
$ws.Range('E11').Value = '  -1.15%  '
$ws.Range('D12').Value = $q + '6.434'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').Value = $q + '1.001'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').Value = $q + '20.40'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = $q + '1.805.84'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = $q + '7.333'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').Value = $q + '0.00001142'
$ws.Range('E17').Value = '  +4.08%  '
$ws.Range('D18').Value = $q + '92.28'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = $q + '0.06589'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').Value = $q + '17.24'
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('D22').Value = $q + '6.042'
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('D23').Value = $q + '28.105.12'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('D24').Value = $q + '11.11'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').Value = $q + '2.219'
$ws.Range('E25').Value = '  -1.49%  '
$ws.Range('D26').Value = $q + '157.15'
$ws.Range('E26').Value = '  -2.16%  '
$ws.Range('D27').Value = $q + '2.423'
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = $q + '20.46'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = $q + '2.009.58'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = $q + '128.49'
$ws.Range('E30').Value = '  +4.19%  '
$ws.Range('D31').Value = $q + '0.1088'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').Value = $q + '1.052'
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('D33').Value = $q + '5.618'
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('D34').Value = $q + '3.624'
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').Value = $q + '0.06877'
$ws.Range('E35').Value = '  -5.69%  '
$ws.Range('D36').Value = $q + '9.060'
$ws.Range('E36').Value = '  +4.17%  '
$ws.Range('D37').Value = $q + '0.02323'
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = $q + '0.2166'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').Value = $q + '11.53'
$ws.Range('E39').Value = '  -7.58%  '
$ws.Range('D40').Value = $q + '5.007'
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('D41').Value = $q + '0.6157'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = $q + '1.153'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('D44').Value = $q + '13.26'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').Value = $q + '0.5938'
$ws.Range('E45').Value = '  -1.86%  '
$ws.Range('D46').Value = $q + '1.289'
$ws.Range('E46').Value = '  -6.08%  '
$ws.Range('E47').Value = '  -1.88%  '
$ws.Range('D48').Value = $q + '124.91'
$ws.Range('E48').Value = '  -1.75%  '
$ws.Range('D49').Value = $q + '1.960'
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('D50').Value = $q + '1.183'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('D51').Value = $q + '0.06762'
$ws.Range('E51').Value = '  -0.87%  '
